{"js": "const body = context.document.body;\n{\n  const results = body.search(\"49\u00d719=931\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '49\u00d719=931' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"56\u00d732=1792\", \"Replace\");\n}\n{\n  const results = body.search(\"54\u00d791=4914\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '54\u00d791=4914' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"16\u00d798=1568\", \"Replace\");\n}\n{\n  const results = body.search(\"83\u00d790=7470\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '83\u00d790=7470' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"65\u00d722=1430\", \"Replace\");\n}\n{\n  const results = body.search(\"72\u00d731=2232\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '72\u00d731=2232' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"31\u00d779=2449\", \"Replace\");\n}\n{\n  const results = body.search(\"11\u00d760=660\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '11\u00d760=660' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"85\u00d797=8245\", \"Replace\");\n}\n{\n  const results = body.search(\"88\u00d762=5456\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '88\u00d762=5456' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"55\u00d747=2585\", \"Replace\");\n}\n{\n  const results = body.search(\"71\u00d775=5325\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '71\u00d775=5325' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"84\u00d743=3612\", \"Replace\");\n}\n{\n  const results = body.search(\"85\u00d784=7140\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '85\u00d784=7140' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"45\u00d711=495\", \"Replace\");\n}\n{\n  const results = body.search(\"65\u00d736=2340\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '65\u00d736=2340' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"88\u00d734=2992\", \"Replace\");\n}\n{\n  const results = body.search(\"97\u00d778=7566\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '97\u00d778=7566' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"30\u00d737=1110\", \"Replace\");\n}\n{\n  const results = body.search(\"51\u00d713=663\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '51\u00d713=663' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"37\u00d768=2516\", \"Replace\");\n}\n{\n  const results = body.search(\"64\u00d787=5568\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '64\u00d787=5568' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"23\u00d717=391\", \"Replace\");\n}\n{\n  const results = body.search(\"55\u00d771=3905\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '55\u00d771=3905' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"15\u00d750=750\", \"Replace\");\n}\n{\n  const results = body.search(\"70\u00d732=2240\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '70\u00d732=2240' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"53\u00d783=4399\", \"Replace\");\n}\n{\n  const results = body.search(\"52\u00d770=3640\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '52\u00d770=3640' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"37\u00d781=2997\", \"Replace\");\n}\n{\n  const results = body.search(\"53\u00d792=4876\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '53\u00d792=4876' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"49\u00d754=2646\", \"Replace\");\n}\n{\n  const results = body.search(\"50\u00d799=4950\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '50\u00d799=4950' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"28\u00d737=1036\", \"Replace\");\n}\n{\n  const results = body.search(\"43\u00d777=3311\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '43\u00d777=3311' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"91\u00d752=4732\", \"Replace\");\n}\n{\n  const results = body.search(\"20\u00d764=1280\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '20\u00d764=1280' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"35\u00d752=1820\", \"Replace\");\n}\n{\n  const results = body.search(\"55\u00d719=1045\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '55\u00d719=1045' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"96\u00d726=2496\", \"Replace\");\n}\n{\n  const results = body.search(\"50\u00d746=2300\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '50\u00d746=2300' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"58\u00d757=3306\", \"Replace\");\n}\n{\n  const results = body.search(\"87\u00d760=5220\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '87\u00d760=5220' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"88\u00d755=4840\", \"Replace\");\n}\n{\n  const results = body.search(\"39\u00d736=1404\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '39\u00d736=1404' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"39\u00d787=3393\", \"Replace\");\n}\n{\n  const results = body.search(\"15\u00d787=1305\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '15\u00d787=1305' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"20\u00d732=640\", \"Replace\");\n}\n{\n  const results = body.search(\"72\u00d748=3456\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected 1 match for '72\u00d748=3456' but found \" + results.items.length);\n  }\n  results.items[0].insertText(\"66\u00d721=1386\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$d.Content.Find.Execute(\"49\u00d719=931\", $false, $false, $false, $false, $false, $true, 1, $false, \"56\u00d732=1792\", 2) | Out-Null\n$d.Content.Find.Execute(\"54\u00d791=4914\", $false, $false, $false, $false, $false, $true, 1, $false, \"16\u00d798=1568\", 2) | Out-Null\n$d.Content.Find.Execute(\"83\u00d790=7470\", $false, $false, $false, $false, $false, $true, 1, $false, \"65\u00d722=1430\", 2) | Out-Null\n$d.Content.Find.Execute(\"72\u00d731=2232\", $false, $false, $false, $false, $false, $true, 1, $false, \"31\u00d779=2449\", 2) | Out-Null\n$d.Content.Find.Execute(\"11\u00d760=660\", $false, $false, $false, $false, $false, $true, 1, $false, \"85\u00d797=8245\", 2) | Out-Null\n$d.Content.Find.Execute(\"88\u00d762=5456\", $false, $false, $false, $false, $false, $true, 1, $false, \"55\u00d747=2585\", 2) | Out-Null\n$d.Content.Find.Execute(\"71\u00d775=5325\", $false, $false, $false, $false, $false, $true, 1, $false, \"84\u00d743=3612\", 2) | Out-Null\n$d.Content.Find.Execute(\"85\u00d784=7140\", $false, $false, $false, $false, $false, $true, 1, $false, \"45\u00d711=495\", 2) | Out-Null\n$d.Content.Find.Execute(\"65\u00d736=2340\", $false, $false, $false, $false, $false, $true, 1, $false, \"88\u00d734=2992\", 2) | Out-Null\n$d.Content.Find.Execute(\"97\u00d778=7566\", $false, $false, $false, $false, $false, $true, 1, $false, \"30\u00d737=1110\", 2) | Out-Null\n$d.Content.Find.Execute(\"51\u00d713=663\", $false, $false, $false, $false, $false, $true, 1, $false, \"37\u00d768=2516\", 2) | Out-Null\n$d.Content.Find.Execute(\"64\u00d787=5568\", $false, $false, $false, $false, $false, $true, 1, $false, \"23\u00d717=391\", 2) | Out-Null\n$d.Content.Find.Execute(\"55\u00d771=3905\", $false, $false, $false, $false, $false, $true, 1, $false, \"15\u00d750=750\", 2) | Out-Null\n$d.Content.Find.Execute(\"70\u00d732=2240\", $false, $false, $false, $false, $false, $true, 1, $false, \"53\u00d783=4399\", 2) | Out-Null\n$d.Content.Find.Execute(\"52\u00d770=3640\", $false, $false, $false, $false, $false, $true, 1, $false, \"37\u00d781=2997\", 2) | Out-Null\n$d.Content.Find.Execute(\"53\u00d792=4876\", $false, $false, $false, $false, $false, $true, 1, $false, \"49\u00d754=2646\", 2) | Out-Null\n$d.Content.Find.Execute(\"50\u00d799=4950\", $false, $false, $false, $false, $false, $true, 1, $false, \"28\u00d737=1036\", 2) | Out-Null\n$d.Content.Find.Execute(\"43\u00d777=3311\", $false, $false, $false, $false, $false, $true, 1, $false, \"91\u00d752=4732\", 2) | Out-Null\n$d.Content.Find.Execute(\"20\u00d764=1280\", $false, $false, $false, $false, $false, $true, 1, $false, \"35\u00d752=1820\", 2) | Out-Null\n$d.Content.Find.Execute(\"55\u00d719=1045\", $false, $false, $false, $false, $false, $true, 1, $false, \"96\u00d726=2496\", 2) | Out-Null\n$d.Content.Find.Execute(\"50\u00d746=2300\", $false, $false, $false, $false, $false, $true, 1, $false, \"58\u00d757=3306\", 2) | Out-Null\n$d.Content.Find.Execute(\"87\u00d760=5220\", $false, $false, $false, $false, $false, $true, 1, $false, \"88\u00d755=4840\", 2) | Out-Null\n$d.Content.Find.Execute(\"39\u00d736=1404\", $false, $false, $false, $false, $false, $true, 1, $false, \"39\u00d787=3393\", 2) | Out-Null\n$d.Content.Find.Execute(\"15\u00d787=1305\", $false, $false, $false, $false, $false, $true, 1, $false, \"20\u00d732=640\", 2) | Out-Null\n$d.Content.Find.Execute(\"72\u00d748=3456\", $false, $false, $false, $false, $false, $true, 1, $false, \"66\u00d721=1386\", 2) | Out-Null\n"}
